$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2:F3 and the "Status" column (C) on the zh-cn / de-de sheets)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Generate-report-for-handback: fill in "Latest Target File",
#    "Latest Handback File" and "Latest Handback DateTime" for each locale,
#    and hyperlink the new "Latest Target File" cell to the source doc (same
#    link used by column A for that row).
# ---------------------------------------------------------------------------

# -- zh-cn --------------------------------------------------------------
$zhcn.Range("I2").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.md"
$zhcn.Range("J2").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.93c95854e16b1616a0761d7d9acba8bd20766fdf.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-07 05:49:33"

$zhcn.Range("I3").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.md"
$zhcn.Range("J3").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.4124c72963f29ffa7603a473b77e015cf4fd79dc.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-07 05:49:33"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md", "", "", "8d52b95f-5e77-4b29-9701-277e4ab73e11.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/a474ad6f-392e-4361-909c-9ce03469b8ae.md", "", "", "a474ad6f-392e-4361-909c-9ce03469b8ae.md")

# -- de-de --------------------------------------------------------------
$dede.Range("I2").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.md"
$dede.Range("J2").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.93c95854e16b1616a0761d7d9acba8bd20766fdf.de-de.xlf"
$dede.Range("K2").Value = "2016-09-07 05:49:51"

$dede.Range("I3").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.md"
$dede.Range("J3").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.4124c72963f29ffa7603a473b77e015cf4fd79dc.de-de.xlf"
$dede.Range("K3").Value = "2016-09-07 05:49:51"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md", "", "", "8d52b95f-5e77-4b29-9701-277e4ab73e11.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/a474ad6f-392e-4361-909c-9ce03469b8ae.md", "", "", "a474ad6f-392e-4361-909c-9ce03469b8ae.md")

# ---------------------------------------------------------------------------
# 3. Widen the columns that now hold the longer status text / file names.
# ---------------------------------------------------------------------------
$ovw.Columns("E").ColumnWidth = 29.9777050018311
$ovw.Columns("F").ColumnWidth = 29.9777050018311

foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns("C").ColumnWidth = 29.9777050018311
    $ws.Columns("I").ColumnWidth = 40
    $ws.Columns("J").ColumnWidth = 40
}

Write-Host "Generated handback report"
